$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.094.90"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "1.910.17"
$ws.Range("E3").Value = "  +1.74%  "
$ws.Range("D4").Value = "'1.009"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.71%  "
$ws.Range("D5").Value = "'316.67"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("E6").Value = "  -0.79%  "
$ws.Range("D7").Value = "'0.4818"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "'0.3815"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.60%  "
$ws.Range("D9").Value = "'0.07357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.50%  "
$ws.Range("D10").Value = "'0.9338"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("D11").Value = "'20.74"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.09%  "
$ws.Range("D12").Value = "'0.07825"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "1.895.87"
$ws.Range("E13").Value = "  +0.69%  "
$ws.Range("D14").Value = "'5.492"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.60%  "
$ws.Range("D15").Value = "'6.620"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.15%  "
$ws.Range("D16").Value = "'91.19"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.13%  "
$ws.Range("D17").Value = "'1.010"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.59%  "
$ws.Range("D18").Value = "'0.000008816"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.02%  "
$ws.Range("D19").Value = "'1.008"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.58%  "
$ws.Range("D20").Value = "28.113.61"
$ws.Range("E20").Value = "  +1.84%  "
$ws.Range("D21").Value = "'14.81"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").Value = "'5.151"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("D23").Value = "2.154.99"
$ws.Range("E23").Value = "  +1.88%  "
$ws.Range("D24").Value = "'10.86"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").Value = "'156.68"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.89%  "
$ws.Range("D26").Value = "'1.914"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.77%  "
$ws.Range("D27").Value = "'18.54"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.33%  "
$ws.Range("D28").Value = "'2.088"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.99%  "
$ws.Range("D29").Value = "'116.26"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.15%  "
$ws.Range("D30").Value = "'4.950"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.40%  "
$ws.Range("D31").Value = "'0.08902"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").Value = "'3.369"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.24%  "
$ws.Range("D33").Value = "'1.243"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.29%  "
$ws.Range("D34").Value = "'0.7661"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("D35").Value = "'4.678"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.09%  "
$ws.Range("D36").Value = "'2.593"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.84%  "
$ws.Range("D37").Value = "'0.02042"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.83%  "
$ws.Range("D38").Value = "'1.097"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.34%  "
$ws.Range("D39").Value = "'0.05293"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.43%  "
$ws.Range("D40").Value = "'0.5490"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").Value = "'2.997"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("D42").Value = "'7.023"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.40%  "
$ws.Range("D43").Value = "'0.1522"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.35%  "
$ws.Range("D44").Value = "'8.422"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.36%  "
$ws.Range("D45").Value = "'10.72"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.25%  "
$ws.Range("D46").Value = "'0.4835"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "
$ws.Range("D47").Value = "'106.78"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.43%  "
$ws.Range("E48").Value = "  -0.73%  "
$ws.Range("D49").Value = "'1.655"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'68.47"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("D51").Value = "'0.06099"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.28%  "
